$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback"
#
# For both the "zh-cn" and "de-de" status sheets, row 6 (the
# 4acf2330-e221-469f-addf-7ac6ca0a9c7c.md entry) now has a handback on file,
# but it is not the most recent handback revision, so the report records the
# stale handback info (Latest Target File / Latest Handback File / Latest
# Handback DateTime) plus a warning in the Error Detail column. The Error
# Detail column (P) is also widened so the message is readable.
# ---------------------------------------------------------------------------

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58c45f83ab50fac49a7048a5d9502d4ce1023215/e2e/4acf2330-e221-469f-addf-7ac6ca0a9c7c.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb0ba353e1fa22888546b8c650db8787d908417b/e2e/4acf2330-e221-469f-addf-7ac6ca0a9c7c.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
$mdName = "4acf2330-e221-469f-addf-7ac6ca0a9c7c.md"

function Update-StatusSheet {
    param(
        $ws,
        [string]$xlfName,
        [string]$handbackDateTime
    )

    # Widen the "Error Detail" column (P / column 16) so the message fits.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I6: Latest Target File -> the handback markdown file name, hyperlinked
    # to its (outdated / "current") source revision, styled like the other
    # hyperlink cells in the sheet (underlined, hyperlink blue).
    $ws.Range("I6").Value = $mdName
    $ws.Range("I6").Style = "HyperLink"
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I6"), $currentUrl, "", "", $mdName)

    # J6: Latest Handback File -> the localized xlf that was handed back.
    $ws.Range("J6").Value = $xlfName

    # K6: Latest Handback DateTime -> when that handback came in.
    $ws.Range("K6").Value = $handbackDateTime

    # P6: Error Detail -> explain that the handback isn't the latest rev.
    $ws.Range("P6").Value = $errorDetail
}

$zh = $wb.Worksheets.Item("zh-cn")
Update-StatusSheet -ws $zh -xlfName "4acf2330-e221-469f-addf-7ac6ca0a9c7c.855bbf43e0358e38371a11dab2d3bbef3f2c276e.zh-cn.xlf" -handbackDateTime "2016-08-23 22:40:46"

$de = $wb.Worksheets.Item("de-de")
Update-StatusSheet -ws $de -xlfName "4acf2330-e221-469f-addf-7ac6ca0a9c7c.855bbf43e0358e38371a11dab2d3bbef3f2c276e.de-de.xlf" -handbackDateTime "2016-08-23 22:40:53"
